$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.78165024612718
$ws.Range("C2").Value = 0.0462677206536974
$ws.Range("D2").Value = 2.60110468262649
$ws.Range("E2").Value = 17.1927294033082
$ws.Range("F2").Value = 0.293155914948531
$ws.Range("B3").Value = 1.48752317521006
$ws.Range("C3").Value = 0.0931246602656657
$ws.Range("D3").Value = 2.27455856484534
$ws.Range("E3").Value = 17.154482281475
$ws.Range("F3").Value = 0.294606855997772
$ws.Range("B4").Value = 1.52782081889351
$ws.Range("C4").Value = 0.0189271013038238
$ws.Range("D4").Value = 2.19004186014249
$ws.Range("E4").Value = 10.8057079897784
$ws.Range("F4").Value = 0.280770388732406
$ws.Range("B5").Value = 1.53494581371079
$ws.Range("C5").Value = 0.201263516911003
$ws.Range("D5").Value = 3.20295795657688
$ws.Range("E5").Value = 11.3971193866693
$ws.Range("F5").Value = 0.279682858006507
$ws.Range("B6").Value = 1.70532292172167
$ws.Range("C6").Value = 0.221873265784416
$ws.Range("D6").Value = 3.05104709318748
$ws.Range("E6").Value = 9.76605881395712
$ws.Range("F6").Value = 0.315273133449451
$ws.Range("B7").Value = 0.237663747213394
$ws.Range("C7").Value = 0.0941478375844147
$ws.Range("D7").Value = 2.55789572688989
$ws.Range("E7").Value = 14.7153654763489
$ws.Range("F7").Value = 0.252772528117768
$ws.Range("B8").Value = 1.3242874451206
$ws.Range("C8").Value = 0.114949280410448
$ws.Range("D8").Value = 3.57332121508243
$ws.Range("E8").Value = 19.6479060346112
$ws.Range("F8").Value = 0.283505694124434
$ws.Range("B9").Value = 1.48797966376133
$ws.Range("C9").Value = 0.12169906534515
$ws.Range("D9").Value = 3.15047089384601
$ws.Range("E9").Value = 17.6944363195715
$ws.Range("F9").Value = 0.191112747764904
$ws.Range("B10").Value = 2.49061555204243
$ws.Range("C10").Value = 0.326108184536942
$ws.Range("D10").Value = 2.2880027503096
$ws.Range("E10").Value = 6.09242291821957
$ws.Range("F10").Value = 0.264849241547033
$ws.Range("B11").Value = 1.9655893232083
$ws.Range("C11").Value = 0.0079330087685354
$ws.Range("D11").Value = 9.72307103536856
$ws.Range("E11").Value = 1.10621691308703
$ws.Range("F11").Value = 0.34600489501338
$ws.Range("B12").Value = 2.87236686642696
$ws.Range("C12").Value = 0.148793028872854
$ws.Range("D12").Value = 9.53129954411118
$ws.Range("E12").Value = 1.3387948104476
$ws.Range("F12").Value = 0.310318611519273
$ws.Range("B13").Value = 2.38444881225538
$ws.Range("C13").Value = 0.102663215220154
$ws.Range("D13").Value = 8.89167039632605
$ws.Range("E13").Value = 1.00163791670994
$ws.Range("F13").Value = 0.209450183761338
$ws.Range("B14").Value = 2.92211898257521
$ws.Range("C14").Value = 0.0328264395717897
$ws.Range("D14").Value = 7.93633312055611
$ws.Range("E14").Value = 1.08461055646884
$ws.Range("F14").Value = 0.504732375831644
$ws.Range("B15").Value = 2.64717443181663
$ws.Range("C15").Value = 0.00479278634955926
$ws.Range("D15").Value = 9.4199822721863
$ws.Range("E15").Value = 1.10336531956178
$ws.Range("F15").Value = 0.40472710482192
$ws.Range("B16").Value = 2.64447195832954
$ws.Range("C16").Value = 0.0422176131885986
$ws.Range("D16").Value = 9.86028746767746
$ws.Range("E16").Value = 1.06888829472491
$ws.Range("F16").Value = 0.364276564743276
$ws.Range("B17").Value = 2.97342186252037
$ws.Range("C17").Value = 0.107225090667447
$ws.Range("D17").Value = 9.89229683168545
$ws.Range("E17").Value = 1.1421399316246
$ws.Range("F17").Value = 0.491225837711439
$ws.Range("B18").Value = 2.52669233159802
$ws.Range("C18").Value = 0.0352239928369997
$ws.Range("D18").Value = 9.51233665315446
$ws.Range("E18").Value = 1.13760439225834
$ws.Range("F18").Value = 0.66884164651133
$ws.Range("B19").Value = 2.59891971683159
$ws.Range("C19").Value = 0.0755602853609649
$ws.Range("D19").Value = 9.84818485082192
$ws.Range("E19").Value = 1.14457050426866
$ws.Range("F19").Value = 0.426105500416268
$ws.Range("B20").Value = 2.71971288037191
$ws.Range("C20").Value = 0.060844263830471
$ws.Range("D20").Value = 1.94816286356057
$ws.Range("E20").Value = 1.3562274243774
$ws.Range("F20").Value = 0.22139506139315
$ws.Range("B21").Value = 2.70024512622127
$ws.Range("C21").Value = 0.00152135123693843
$ws.Range("D21").Value = 2.10544418350467
$ws.Range("E21").Value = 1.35463463721079
$ws.Range("F21").Value = 0.144753565272957
$ws.Range("B22").Value = 2.61615332507346
$ws.Range("C22").Value = 0.14231616399217
$ws.Range("D22").Value = 1.92286758152138
$ws.Range("E22").Value = 1.33112848464231
$ws.Range("F22").Value = 0.36925342229849
$ws.Range("B23").Value = 2.99126087542907
$ws.Range("C23").Value = 0.122618806023828
$ws.Range("D23").Value = 1.82097858022822
$ws.Range("E23").Value = 1.88165917698091
$ws.Range("F23").Value = 0.541713404557175
$ws.Range("B24").Value = 2.71821010723278
$ws.Range("C24").Value = 0.0453390989277372
$ws.Range("D24").Value = 2.06034365650971
$ws.Range("E24").Value = 1.83382826454228
$ws.Range("F24").Value = 0.531434015645634
$ws.Range("B25").Value = 2.55294073929473
$ws.Range("C25").Value = 0.0480180760217598
$ws.Range("D25").Value = 1.49604663168749
$ws.Range("E25").Value = 1.32628253315006
$ws.Range("F25").Value = 0.496835763663412
$ws.Range("B26").Value = 2.25228031617351
$ws.Range("C26").Value = 0.0471863303992923
$ws.Range("D26").Value = 2.09483419799816
$ws.Range("E26").Value = 1.2749897567587
$ws.Range("F26").Value = 0.808261193540028
$ws.Range("B27").Value = 2.67378090274662
$ws.Range("C27").Value = 0.0798211773715247
$ws.Range("D27").Value = 2.06503113422661
$ws.Range("E27").Value = 1.56187570474451
$ws.Range("F27").Value = 0.810943107927034
$ws.Range("B28").Value = 2.95831085724958
$ws.Range("C28").Value = 0.0208776940691039
$ws.Range("D28").Value = 1.8609369028358
$ws.Range("E28").Value = 2.38878798913279
$ws.Range("F28").Value = 0.669812599273033

$ws.Range("I23").Select()
